$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: quarter-period headers (shift window forward by one quarter)
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# Row 9: publish-date headers (shift window forward by one quarter)
$ws.Range("D9").Value = "1400-10-29 (2)"
$ws.Range("E9").Value = "1401-04-08 (8)"
$ws.Range("F9").Value = "1401-05-04 (3)"
$ws.Range("G9").Value = "1401-08-30 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-30 (8)"
$ws.Range("J9").Value = "1401-05-04 (2)"
$ws.Range("K9").Value = "1401-08-30 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-30"

# Financial data rows: shift one quarter + recompute (read_price change for row 12/13/14/17)
# Row 11
$ws.Range("D11").Value = 211182
$ws.Range("E11").Value = 162845
$ws.Range("F11").Value = 234396
$ws.Range("G11").Value = 207527
$ws.Range("H11").Value = 159961
$ws.Range("I11").Value = 138812
$ws.Range("J11").Value = 236542
$ws.Range("K11").Value = 187229
$ws.Range("L11").Value = 133430
$ws.Range("M11").Value = 88615

# Row 12
$ws.Range("D12").Value = -80536
$ws.Range("E12").Value = -82861
$ws.Range("F12").Value = -102279
$ws.Range("G12").Value = -145235
$ws.Range("H12").Value = -181308
$ws.Range("I12").Value = -119285
$ws.Range("J12").Value = -150219
$ws.Range("K12").Value = -145770
$ws.Range("L12").Value = -86446
$ws.Range("M12").Value = -61521

# Row 13
$ws.Range("D13").Value = 130646
$ws.Range("E13").Value = 79983
$ws.Range("F13").Value = 132116
$ws.Range("G13").Value = 62292
$ws.Range("H13").Value = -21347
$ws.Range("I13").Value = 19527
$ws.Range("J13").Value = 86323
$ws.Range("K13").Value = 41459
$ws.Range("L13").Value = 46983
$ws.Range("M13").Value = 27094

# Row 14
$ws.Range("D14").Value = -42382
$ws.Range("E14").Value = -38281
$ws.Range("F14").Value = -37942
$ws.Range("G14").Value = -40861
$ws.Range("H14").Value = -16717
$ws.Range("I14").Value = -31010
$ws.Range("J14").Value = -43652
$ws.Range("K14").Value = -44386
$ws.Range("L14").Value = -24090
$ws.Range("M14").Value = -24943

# Row 16
$ws.Range("D16").Value = 9928
$ws.Range("E16").Value = -19989
$ws.Range("F16").Value = -9081
$ws.Range("G16").Value = 7900
$ws.Range("H16").Value = 733
$ws.Range("I16").Value = -3055
$ws.Range("J16").Value = 3306
$ws.Range("K16").Value = -2309
$ws.Range("L16").Value = 13656
$ws.Range("M16").Value = 12849

# Row 17
$ws.Range("D17").Value = 98192
$ws.Range("E17").Value = 21713
$ws.Range("F17").Value = 85093
$ws.Range("G17").Value = 29331
$ws.Range("H17").Value = -37331
$ws.Range("I17").Value = -14539
$ws.Range("J17").Value = 45976
$ws.Range("K17").Value = -5236
$ws.Range("L17").Value = 36549
$ws.Range("M17").Value = 15000

# Row 19
$ws.Range("D19").Value = -929
$ws.Range("E19").Value = 907
$ws.Range("F19").Value = 3165
$ws.Range("G19").Value = 4186
$ws.Range("H19").Value = 2694
$ws.Range("I19").Value = 4137
$ws.Range("J19").Value = 699
$ws.Range("K19").Value = 4990
$ws.Range("L19").Value = 1677
$ws.Range("M19").Value = -4653

# Row 20
$ws.Range("D20").Value = 97263
$ws.Range("E20").Value = 22620
$ws.Range("F20").Value = 88258
$ws.Range("G20").Value = 33517
$ws.Range("H20").Value = -34637
$ws.Range("I20").Value = -10401
$ws.Range("J20").Value = 46674
$ws.Range("K20").Value = -246
$ws.Range("L20").Value = 38225
$ws.Range("M20").Value = 10347

# Row 21
$ws.Range("J21").Value = -9335
$ws.Range("K21").Value = 8844
$ws.Range("L21").Value = -15496
$ws.Range("M21").Value = 11845

# Row 22
$ws.Range("D22").Value = 97263
$ws.Range("E22").Value = 22620
$ws.Range("F22").Value = 88258
$ws.Range("G22").Value = 33517
$ws.Range("H22").Value = -34637
$ws.Range("I22").Value = -10401
$ws.Range("J22").Value = 37339
$ws.Range("K22").Value = 8598
$ws.Range("L22").Value = 22729
$ws.Range("M22").Value = 22192

# Row 24
$ws.Range("D24").Value = 97263
$ws.Range("E24").Value = 22620
$ws.Range("F24").Value = 88258
$ws.Range("G24").Value = 33517
$ws.Range("H24").Value = -34637
$ws.Range("I24").Value = -10401
$ws.Range("J24").Value = 37339
$ws.Range("K24").Value = 8598
$ws.Range("L24").Value = 22729
$ws.Range("M24").Value = 22192

# Row 26
$ws.Range("D26").Value = 8831
$ws.Range("E26").Value = 9826
$ws.Range("F26").Value = 10287
$ws.Range("G26").Value = 9195
$ws.Range("H26").Value = 8429
$ws.Range("I26").Value = 8713
$ws.Range("J26").Value = 8168
$ws.Range("K26").Value = 7738
$ws.Range("L26").Value = 6908
$ws.Range("M26").Value = 5280

